# Apply the updated crypto price/volume snapshot values (and the
# Polkadot/WrappedEther row swap) to worksheet 1, matching the new
# GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.580.12"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "'1.924.07"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "'326.46"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'0.4057"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").Value = "'0.08221"
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("D11").Value = "'23.89"
$ws.Range("E11").Value = "  +1.74%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.910.44"
$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.122"
$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").Value = "'7.297"
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("D15").Value = "'91.75"
$ws.Range("E15").Value = "  +1.70%  "

$ws.Range("D16").Value = "'0.06860"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").Value = "'1.012"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").Value = "'17.65"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D21").Value = "'29.585.32"

$ws.Range("D22").Value = "'5.678"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("D23").Value = "'12.01"
$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").Value = "'2.145.01"
$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("D26").Value = "'156.08"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'6.411"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").Value = "'2.093"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").Value = "'120.73"
$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("D32").Value = "'0.09604"
$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("D33").Value = "'5.615"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'1.380"
$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").Value = "'0.06359"
$ws.Range("E36").Value = "  +4.30%  "

$ws.Range("D37").Value = "'0.02291"
$ws.Range("E37").Value = "  +0.99%  "

$ws.Range("D38").Value = "'1.191"
$ws.Range("E38").Value = "  +1.35%  "

$ws.Range("D39").Value = "'0.5958"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("D40").Value = "'10.75"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D42").Value = "'7.873"
$ws.Range("E42").Value = "  -1.45%  "

$ws.Range("D43").Value = "'0.1852"
$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("E44").Value = "  +1.09%  "

$ws.Range("D45").Value = "'1.249"
$ws.Range("E45").Value = "  -2.72%  "

$ws.Range("D46").Value = "'12.45"
$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("D47").Value = "'0.07542"
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("D48").Value = "'0.5561"

$ws.Range("D49").Value = "'1.988"
$ws.Range("E49").Value = "  +2.46%  "

$ws.Range("D50").Value = "'119.35"
$ws.Range("E50").Value = "  +3.03%  "

$ws.Range("D51").Value = "'2.434"
$ws.Range("E51").Value = "  +0.73%  "
